$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "37.723.58"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "  -0.09%  "
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.080.19"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "  -1.67%  "
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "  +0.01%  "
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "234.14"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "  -0.68%  "
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "  -0.71%  "
$c.Style = "Normal"
$c = $ws.Range("B7")
$c.NumberFormat = "@"
$c.Value = "USDC"
$c.Style = "Normal"
$c = $ws.Range("C7")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "  -0.04%  "
$c.Style = "Normal"
$c = $ws.Range("B8")
$c.NumberFormat = "@"
$c.Value = "Solana"
$c.Style = "Normal"
$c = $ws.Range("C8")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "58.33"
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "  -0.20%  "
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "  -0.36%  "
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "  +0.06%  "
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "  +2.62%  "
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "14.97"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "  +2.58%  "
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "2.386.32"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "  -1.55%  "
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "20.97"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "  -1.72%  "
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.772"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "  -2.28%  "
$c.Style = "Normal"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "5.30"
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "  +1.26%  "
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "2.090.63"
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "  -1.03%  "
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "37.625.30"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "  -0.30%  "
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.17"
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "  -0.62%  "
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "71.14"
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "  +1.13%  "
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "  +0.84%  "
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "228.07"
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "  +0.03%  "
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "  -0.13%  "
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "  -0.29%  "
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.39"
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "  -1.54%  "
$c.Style = "Normal"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "168.97"
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "  +0.18%  "
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "  +2.36%  "
$c.Style = "Normal"
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = "  -0.40%  "
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = "  -1.51%  "
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "19.45"
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = "  -0.20%  "
$c.Style = "Normal"
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = "  +1.78%  "
$c.Style = "Normal"
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = "  +0.43%  "
$c.Style = "Normal"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.0630"
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = "  +1.10%  "
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = "  +1.21%  "
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = "  -4.04%  "
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = "  +2.62%  "
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = "  -4.29%  "
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = "  -0.03%  "
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "  -5.28%  "
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.0976"
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "  +1.01%  "
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "98.11"
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "  +0.23%  "
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "  +0.41%  "
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "1.450.81"
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "  -1.95%  "
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "16.61"
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "  +6.20%  "
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "4.29"
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "  +1.52%  "
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "  -0.64%  "
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "  +0.62%  "
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "7.38"
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "  +0.63%  "
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "  -0.40%  "
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "2.270.44"
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "  -1.71%  "
$c.Style = "Normal"
